$d = $word.ActiveDocument
$d.Content.Find.Execute("CUALQUIER HIJO.", $true, $false, $false, $false, $false, $true, 1, $false, "CUALQUIER HIJO." + [char]9 + "MULTIPLE HIJOS", 2)
